$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 11360
$ws.Range("F4").Value = 1286
$ws.Range("F5").Value = 1154
$ws.Range("F7").Value = 1217
$ws.Range("F9").Value = 947
$ws.Range("F11").Value = 2260
$ws.Range("F13").Value = 1109
$ws.Range("F14").Value = 868
$ws.Range("F16").Value = 854
$ws.Range("F17").Value = 1005
$ws.Range("F19").Value = 101
$ws.Range("F20").Value = 682
$ws.Range("F21").Value = 710
$ws.Range("F22").Value = 151
$ws.Range("F23").Value = 397
$ws.Range("F24").Value = 1054
$ws.Range("F25").Value = 62
$ws.Range("F26").Value = 438
$ws.Range("F27").Value = 533
$ws.Range("F29").Value = 267
$ws.Range("F30").Value = 266
$ws.Range("F31").Value = 631
$ws.Range("F32").Value = 2507
$ws.Range("F33").Value = 432
$ws.Range("F35").Value = 149
$ws.Range("F36").Value = 12
$ws.Range("D37").Value = "沙溪大道沙溪地铁站E2出口桥下 飞梦篮球公园(沙溪店)"
$ws.Range("F37").Value = 67
$ws.Range("F38").Value = 1504
$ws.Range("F40").Value = 132
$ws.Range("F41").Value = 64
$ws.Range("F44").Value = 8
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 156
$ws.Range("F16").Value = 93
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 674
$ws.Range("F4").Value = 636
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 11360
$ws.Range("F4").Value = 636
$ws.Range("F5").Value = 1154
$ws.Range("F8").Value = 1217
$ws.Range("F12").Value = 947
$ws.Range("F13").Value = 2260
$ws.Range("F15").Value = 1109
$ws.Range("F16").Value = 868
$ws.Range("F18").Value = 854
$ws.Range("F19").Value = 1005
$ws.Range("F22").Value = 101
$ws.Range("F23").Value = 682
$ws.Range("F26").Value = 710
$ws.Range("F27").Value = 151
$ws.Range("F28").Value = 397
$ws.Range("F29").Value = 1054
$ws.Range("F30").Value = 62
$ws.Range("F31").Value = 438
$ws.Range("F32").Value = 533
$ws.Range("F34").Value = 267
$ws.Range("F35").Value = 2507
$ws.Range("F36").Value = 156
$ws.Range("F37").Value = 432
$ws.Range("D38").Value = "沙溪大道沙溪地铁站E2出口桥下 飞梦篮球公园(沙溪店)"
$ws.Range("F38").Value = 67
$ws.Range("F39").Value = 1504
$ws.Range("F41").Value = 132
$ws.Range("F42").Value = 64
$ws.Range("F45").Value = 8
